$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45182 = 2023-09-13) for every
# data row (rows 2-527). The commit updates this "last changed" timestamp to
# serial 45184 (2023-09-15) for all rows.
$ws.Range("C2:C527").Value = 45184
